$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.10"
$ws.Range("E2").Value = "'-0.39%"
$ws.Range("D3").Value = "'26.40"
$ws.Range("E3").Value = "'3.52%"
$ws.Range("D4").Value = "'5.160"
$ws.Range("E4").Value = "'0.98%"
$ws.Range("D5").Value = "'0.05605"
$ws.Range("E5").Value = "'0.29%"
$ws.Range("D6").Value = "'6.468"
$ws.Range("E6").Value = "'-0.09%"
$ws.Range("D7").Value = "'0.8190"
$ws.Range("E7").Value = "'0.09%"
$ws.Range("D8").Value = "'0.8274"
$ws.Range("E8").Value = "'-1.67%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1332"
$ws.Range("E9").Value = "'-0.60%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.06931"
$ws.Range("E10").Value = "'-0.25%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02887"
$ws.Range("E11").Value = "'1.10%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09383"
$ws.Range("E12").Value = "'0.10%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001514"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0005965"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("D15").Value = "'0.006164"
$ws.Range("E15").Value = "'-0.97%"
$ws.Range("D16").Value = "'3.655"
$ws.Range("E16").Value = "'3.74%"
$ws.Range("E17").Value = "'0.58%"
$ws.Range("D18").Value = "'2.182"
$ws.Range("E18").Value = "'4.37%"
$ws.Range("D20").Value = "'0.03078"
$ws.Range("E20").Value = "'-4.26%"
$ws.Range("E21").Value = "'-2.27%"
$ws.Range("D22").Value = "'3.740"
$ws.Range("E22").Value = "'-0.06%"
$ws.Range("D23").Value = "'0.04619"
$ws.Range("E23").Value = "'-1.73%"
$ws.Range("E24").Value = "'-2.46%"
$ws.Range("D25").Value = "'0.001245"
$ws.Range("E25").Value = "'-0.14%"
$ws.Range("D26").Value = "'0.004493"
$ws.Range("E26").Value = "'-2.84%"
$ws.Range("D27").Value = "'0.00009596"
$ws.Range("E27").Value = "'-1.08%"
$ws.Range("D28").Value = "'0.0001936"
$ws.Range("D40").Value = "'0.03648"
$ws.Range("E40").Value = "'-0.40%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1362"
$ws.Range("E41").Value = "'29.63%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006168"
$ws.Range("E42").Value = "'-0.38%"
$ws.Range("D43").Value = "'0.002619"
$ws.Range("E43").Value = "'3.26%"
$ws.Range("D44").Value = "'0.008984"
$ws.Range("E44").Value = "'19.46%"
$ws.Range("D45").Value = "'0.00005353"
$ws.Range("E45").Value = "'0.82%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("E47").Value = "'8.17%"
$ws.Range("D48").Value = "'0.002330"
$ws.Range("E48").Value = "'9.74%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.10%"
